$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: new planned-hours / actual-date / weekday / delta entry ---
$ws.Range("P25").Value = 4
$ws.Range("Q25").Value = (Get-Date -Year 2020 -Month 4 -Day 1).Date
$ws.Range("R25").Value = "Wednesday"
$ws.Range("S25").Formula = "=M25-P25"

# --- Row 26: new planned-hours / actual-date(pending) / weekday / delta entry ---
$ws.Range("P26").Value = 10
$ws.Range("Q26").Value = "2020-04-03 Pending"
$ws.Range("R26").Value = "Friday"
$ws.Range("S26").Formula = "=M26-P26"

# --- Row 29: add delta formula next to the existing actual-date/weekday entry ---
$ws.Range("S29").Formula = "=M29-P29"

# --- Insert a new row before row 42 for the "Guessers Average Component" entry,
#     shifting the summary rows (Sum added tasks / Tot ant tim / Ant dgr / % of plan + added)
#     down by one ---
[void]$ws.Rows("42:42").Insert()

$ws.Range("F42").Value = "Guessers Average Component"
$ws.Range("P42").Value = 4
$ws.Range("Q42").Value = (Get-Date -Year 2020 -Month 4 -Day 2).Date
$ws.Range("R42").Value = "Thursday"

$ws.Range("P43").Formula = "=SUM(P41:P42)"

$ws.Range("P45").Formula = "=P44/7"
$ws.Range("P45").NumberFormat = "0.0"

$ws.Range("M37").Formula = "=M36/7"
$ws.Range("P37").Formula = "=P36/7"
$ws.Range("P37").NumberFormat = "0"

# --- restore the view selection (best effort; frozen-pane scroll position
#     is not independently addressable through this COM surface) ---
[void]$ws.Range("P26").Select()

Write-Output "edit applied"
